# Apply citation-note markup changes described by the commit:
#   "added -nb markup for citation notes"
#
# 1) "...Scientific American –jmr."  ->  "...Scientific American -nb."
# 2) Merge the {Tradition and Design in the Illiad} / p. 222 / %% runs
#    back into a single contiguous run of text (no content change, just
#    collapsing the run split).

$d = $word.ActiveDocument
$enDash = [char]0x2013

# --- 1) Replace the old "-jmr." citation-note marker with "-nb." ---------
$d.Content.Find.Execute(
    "American " + $enDash + "jmr.",   # "American \u2013jmr."
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "American -nb.",
    2
) | Out-Null

# --- 2) Re-join the split "{Tradition and Design in the Illiad} p. 222  %%" ---
$d.Content.Find.Execute(
    "{Tradition and Design in the Illiad} p. 222  %%",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "{Tradition and Design in the Illiad} p. 222  %%",
    2
) | Out-Null
